# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from the last
# existing header cell (AC1) onto the three new header cells.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))

# Set the new header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record for every data row (rows 2-50).
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 76  # AD
    $ws.Cells.Item($r, 31).Value = 86  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
